$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: advance the date in A1 by one day (2024-01-17 -> 2024-01-18)
$ws.Range("A1").Value2 = 45309

# Step 2: update prices
$ws.Range("D30").Value2 = 1475
$ws.Range("D31").Value2 = 1680
